$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-location")

# New postal-code rows for Ben Mansour (BNMR) location, in eng/fra/ara,
# appended after the existing data (rows 110-119).
$rows = @(
    @(110, 10110, 10110, 5, "Postal Code",     "BNMR", "eng"),
    @(111, 10111, 10111, 5, "Postal Code",     "BNMR", "eng"),
    @(112, 10113, 10113, 5, "Postal Code",     "BNMR", "eng"),
    @(113, 10114, 10114, 5, "Postal Code",     "BNMR", "eng"),
    @(114, 10111, 10111, 5, "code postal",     "BNMR", "fra"),
    @(115, 10110, 10110, 5, "code postal",     "BNMR", "fra"),
    @(116, 10113, 10113, 5, "code postal",     "BNMR", "fra"),
    @(117, 10114, 10114, 5, "code postal",     "BNMR", "fra"),
    @(118, 10111, 10111, 5, "الرمز البريدي", "BNMR", "ara"),
    @(119, 10110, 10110, 5, "الرمز البريدي", "BNMR", "ara")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $true
    $ws.Cells.Item($rowNum, 8).Value = "superadmin"
    $ws.Cells.Item($rowNum, 9).Value = "now()"
}

# Clear the stale selection left over on the old "last row" range and
# reset it to the default top-left cell, matching the refreshed sheet view.
$null = $ws.Range("A1").Select()
